# Updates the "Price" (column D) and "Volume(1h)" (column E) columns of the
# cryptocurrency table on the active sheet, mirroring the scraped values
# refreshed by the "Updated cryptos list" GitHub Actions job.
#
# Column D holds price strings that often *look* numeric (e.g. "1.002",
# "165.00", "0.07060"); Excel's COM layer auto-coerces a bare numeric-looking
# string into a real number (losing the trailing zeros / thousands-style
# dots that are meaningful here). Prefixing those values with a leading
# apostrophe forces Excel to store them as literal text, exactly like a user
# typing `'1.002` into a cell - the apostrophe itself is not kept as part of
# the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "31.255.18"; E = "  +2.86%  " },
    @{ Row = 3; D = "2.000.46"; E = "  +6.91%  " },
    @{ Row = 4; D = $null; E = "  +0.26%  " },
    @{ Row = 5; D = "0.7730"; E = "  +63.97%  " },
    @{ Row = 6; D = "258.01"; E = "  +5.48%  " },
    @{ Row = 7; D = "1.002"; E = "  +0.20%  " },
    @{ Row = 8; D = "0.3546"; E = "  +23.64%  " },
    @{ Row = 9; D = "28.83"; E = "  +32.23%  " },
    @{ Row = 10; D = "0.07060"; E = "  +8.78%  " },
    @{ Row = 11; D = "0.8598"; E = "  +18.60%  " },
    @{ Row = 12; D = "0.08194"; E = "  +5.15%  " },
    @{ Row = 13; D = "101.70"; E = "  +1.28%  " },
    @{ Row = 14; D = "2.003.23"; E = "  +7.03%  " },
    @{ Row = 15; D = "5.585"; E = "  +8.04%  " },
    @{ Row = 16; D = "15.48"; E = "  +18.47%  " },
    @{ Row = 17; D = "273.42"; E = "  -4.01%  " },
    @{ Row = 18; D = "31.265.61"; E = "  +2.93%  " },
    @{ Row = 19; D = "5.933"; E = "  +11.30%  " },
    @{ Row = 20; D = "0.000007970"; E = "  +6.43%  " },
    @{ Row = 21; D = "2.269.22"; E = "  +7.38%  " },
    @{ Row = 22; D = "1.002"; E = "  +0.24%  " },
    @{ Row = 23; D = $null; E = "  +0.30%  " },
    @{ Row = 24; D = "7.156"; E = "  +13.41%  " },
    @{ Row = 25; D = "10.05"; E = "  +11.25%  " },
    @{ Row = 26; D = "165.00"; E = "  +1.03%  " },
    @{ Row = 27; D = "0.1476"; E = "  +52.76%  " },
    @{ Row = 28; D = $null; E = "  +5.75%  " },
    @{ Row = 29; D = "2.389"; E = $null },
    @{ Row = 30; D = "1.618"; E = "  +8.82%  " },
    @{ Row = 31; D = "4.625"; E = "  +9.49%  " },
    @{ Row = 32; D = "1.361"; E = "  +3.04%  " },
    @{ Row = 33; D = "4.425"; E = "  +6.84%  " },
    @{ Row = 34; D = "0.05225"; E = $null },
    @{ Row = 35; D = $null; E = "  +8.62%  " },
    @{ Row = 36; D = "0.7737"; E = "  +12.52%  " },
    @{ Row = 37; D = "2.809"; E = "  +3.17%  " },
    @{ Row = 38; D = "0.01998"; E = "  +5.19%  " },
    @{ Row = 39; D = "2.942"; E = "  +3.56%  " },
    @{ Row = 40; D = $null; E = "  +6.95%  " },
    @{ Row = 41; D = "80.04"; E = "  +5.10%  " },
    @{ Row = 42; D = $null; E = "  +12.05%  " },
    @{ Row = 43; D = "2.145"; E = "  +9.32%  " },
    @{ Row = 44; D = "106.79"; E = "  +5.95%  " },
    @{ Row = 45; D = "0.8585"; E = "  +4.32%  " },
    @{ Row = 46; D = "1.003"; E = "  +0.34%  " },
    @{ Row = 47; D = "7.753"; E = "  +10.70%  " },
    @{ Row = 48; D = "9.933"; E = "  +1.51%  " },
    @{ Row = 49; D = "0.4343"; E = "  +11.15%  " },
    @{ Row = 50; D = "36.76"; E = "  +4.96%  " },
    @{ Row = 51; D = "1.516"; E = "  +14.14%  " }
)

# Excel/COM auto-converts a cell's assigned string into a real number
# whenever the string parses cleanly as one (e.g. "1.002", "165.00",
# "0.07060"). These "Price" values must stay literal text (trailing zeros
# and the multi-dot thousands style like "31.255.18" are significant), so
# number-looking values are written with a leading apostrophe - the same
# quote-prefix trick a person uses when typing such a value into a sheet.
# The apostrophe itself is not stored as part of the cell's value/text.
function Set-TextValue($cell, $value) {
    if ($value -match '^-?\d+(\.\d+)?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($row, 4) $u.D
    }
    if ($null -ne $u.E) {
        Set-TextValue $ws.Cells.Item($row, 5) $u.E
    }
}
